# BOM update: the old "3.3V linear regulator" (TPS735, Texas Instruments,
# 6-SON package) line item is replaced with a new "3v linear reg" part
# (TLV1117-33, ti, sot233 package) at the bottom of the first BOM table.
#
# Net effect observed in the OOXML diff: row 16 (old TPS735 row) is removed
# (shifting rows 17-23 up to 16-22, along with everything below, including
# the mirrored/reference table further down the sheet), and a brand new
# row is appended at the end of the first table (now row 23) holding the
# new TLV1117-33 part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 3.3V linear regulator row entirely (shifts everything
# below it up by one row).
$ws.Rows(16).Delete()

# Row-shifting re-stores a few untouched numeric literals with extra
# floating-point noise (e.g. 5.89 -> 5.8899999999999997); write the clean
# values back so those cells round-trip the way they did before the shift.
$ws.Cells.Item(17, 6).Value = 5.89
$ws.Cells.Item(21, 6).Value = 0.51
$ws.Cells.Item(42, 7).Value = 0.17
$ws.Cells.Item(47, 7).Value = 5.89

# Append the new part as a new row after the old row 22 (now the last row
# of the first table, "ferrite bead").
$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = "3v linear reg"
$ws.Cells.Item(23, 4).Value = "ti"
$ws.Cells.Item(23, 5).Value = "https://www.digikey.com/scripts/DkSearch/dksus.dll?Detail&itemSeq=197568672&uq=635991844891259046"
$ws.Cells.Item(23, 7).Value = "sot233"
$ws.Cells.Item(23, 3).Value = "http://www.ti.com/lit/ds/symlink/tlv1117-33.pdf"
$ws.Cells.Item(23, 6).Value = 0.62

# Carry the "Man." column formatting (cellXf s="4") used throughout the
# table onto the new row's D cell (same look as D16:D22 above it).
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(23, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the saved selection / scroll position to match the edited
# workbook (window moved left, selection now sits on the new part's
# datasheet cell).
$ws.Range("C23").Select()
$excel.ActiveWindow.Left = 80
